$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TC's_Sheet")
$ws.Range("A7").Value = 8101
$ws.Range("B7").Value = "ER_W_2s"
$ws.Range("C7").Value = "ER_W_2s"
$ws.Hyperlinks.Add($ws.Cells.Item(7,3), "", "ER_W_2s!A1", "", "ER_W_2s")
